$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 115.2213693333333
$ws.Range("H2").Value = 345.664108
$ws.Range("I2").Value = 0.2787408744545015
$ws.Range("J2").Value = 0.2787408744545015
$ws.Range("M2").Value = 1021.934916333333
$ws.Range("N2").Value = 3065.804749
$ws.Range("O2").Value = 0.8026347959846111
$ws.Range("P2").Value = 0.802634795984611
$ws.Range("Q2").Value = 117748.7404294721
$ws.Range("R2").Value = 1059738.663865249
$ws.Range("S2").Value = 0.2237271249003609
$ws.Range("T2").Value = 0.2237271249003608

# Row 3
$ws.Range("G3").Value = 115.2213693333333
$ws.Range("H3").Value = 345.664108
$ws.Range("I3").Value = 0.2787408744545015
$ws.Range("J3").Value = 0.2787408744545015
$ws.Range("O3").Value = 0.04931810976893385
$ws.Range("P3").Value = 0.04931810976893384
$ws.Range("Q3").Value = 7235.102857122761
$ws.Range("R3").Value = 65115.92571410485
$ws.Range("S3").Value = 0.01374697304343572
$ws.Range("T3").Value = 0.01374697304343571

# Row 4
$ws.Range("G4").Value = 115.2213693333333
$ws.Range("H4").Value = 345.664108
$ws.Range("I4").Value = 0.2787408744545015
$ws.Range("J4").Value = 0.2787408744545015
$ws.Range("M4").Value = 187.139577
$ws.Range("N4").Value = 561.418731
$ws.Range("O4").Value = 0.1469807262726385
$ws.Range("P4").Value = 0.1469807262726385
$ws.Range("Q4").Value = 21562.47831840077
$ws.Range("R4").Value = 194062.3048656069
$ws.Range("S4").Value = 0.04096953616919299
$ws.Range("T4").Value = 0.04096953616919297

# Row 5
$ws.Range("G5").Value = 115.2213693333333
$ws.Range("H5").Value = 345.664108
$ws.Range("I5").Value = 0.2787408744545015
$ws.Range("J5").Value = 0.2787408744545015
$ws.Range("M5").Value = 1.357726666666667
$ws.Range("N5").Value = 4.073180000000001
$ws.Range("O5").Value = 0.001066367973816652
$ws.Range("P5").Value = 0.001066367973816652
$ws.Range("Q5").Value = 156.4391257137156
$ws.Range("R5").Value = 1407.95213142344
$ws.Range("S5").Value = 0.0002972403415119285
$ws.Range("T5").Value = 0.0002972403415119284

# Row 6
$ws.Range("I6").Value = 0.44716501655323
$ws.Range("J6").Value = 0.4471650165532299
$ws.Range("M6").Value = 1021.934916333333
$ws.Range("N6").Value = 3065.804749
$ws.Range("O6").Value = 0.8026347959846111
$ws.Range("P6").Value = 0.802634795984611
$ws.Range("Q6").Value = 188896.2914617725
$ws.Range("R6").Value = 1700066.623155953
$ws.Range("S6").Value = 0.358910201832657
$ws.Range("T6").Value = 0.3589102018326569

# Row 7
$ws.Range("I7").Value = 0.44716501655323
$ws.Range("J7").Value = 0.4471650165532299
$ws.Range("O7").Value = 0.04931810976893385
$ws.Range("P7").Value = 0.04931810976893384
$ws.Range("S7").Value = 0.02205333337119932
$ws.Range("T7").Value = 0.02205333337119931

# Row 8
$ws.Range("I8").Value = 0.44716501655323
$ws.Range("J8").Value = 0.4471650165532299
$ws.Range("M8").Value = 187.139577
$ws.Range("N8").Value = 561.418731
$ws.Range("O8").Value = 0.1469807262726385
$ws.Range("P8").Value = 0.1469807262726385
$ws.Range("Q8").Value = 34591.21663819775
$ws.Range("R8").Value = 311320.9497437798
$ws.Range("S8").Value = 0.06572463889671018
$ws.Range("T8").Value = 0.06572463889671015

# Row 9
$ws.Range("I9").Value = 0.44716501655323
$ws.Range("J9").Value = 0.4471650165532299
$ws.Range("M9").Value = 1.357726666666667
$ws.Range("N9").Value = 4.073180000000001
$ws.Range("O9").Value = 0.001066367973816652
$ws.Range("P9").Value = 0.001066367973816652
$ws.Range("Q9").Value = 250.9646436901201
$ws.Range("R9").Value = 2258.68179321108
$ws.Range("S9").Value = 0.0004768424526635575
$ws.Range("T9").Value = 0.0004768424526635572

# Row 10
$ws.Range("G10").Value = 60.55095666666667
$ws.Range("H10").Value = 181.65287
$ws.Range("I10").Value = 0.1464834753134679
$ws.Range("J10").Value = 0.1464834753134678
$ws.Range("M10").Value = 1021.934916333333
$ws.Range("N10").Value = 3065.804749
$ws.Range("O10").Value = 0.8026347959846111
$ws.Range("P10").Value = 0.802634795984611
$ws.Range("Q10").Value = 61879.13683505329
$ws.Range("R10").Value = 556912.2315154796
$ws.Range("S10").Value = 0.1175727343233421
$ws.Range("T10").Value = 0.1175727343233421

# Row 11
$ws.Range("G11").Value = 60.55095666666667
$ws.Range("H11").Value = 181.65287
$ws.Range("I11").Value = 0.1464834753134679
$ws.Range("J11").Value = 0.1464834753134678
$ws.Range("O11").Value = 0.04931810976893385
$ws.Range("P11").Value = 0.04931810976893384
$ws.Range("Q11").Value = 3802.18011741488
$ws.Range("R11").Value = 34219.62105673391
$ws.Range("S11").Value = 0.007224288114844521
$ws.Range("T11").Value = 0.007224288114844517

# Row 12
$ws.Range("G12").Value = 60.55095666666667
$ws.Range("H12").Value = 181.65287
$ws.Range("I12").Value = 0.1464834753134679
$ws.Range("J12").Value = 0.1464834753134678
$ws.Range("M12").Value = 187.139577
$ws.Range("N12").Value = 561.418731
$ws.Range("O12").Value = 0.1469807262726385
$ws.Range("P12").Value = 0.1469807262726385
$ws.Range("Q12").Value = 11331.48041754533
$ws.Range("R12").Value = 101983.323757908
$ws.Range("S12").Value = 0.02153024758851362
$ws.Range("T12").Value = 0.02153024758851362

# Row 13
$ws.Range("G13").Value = 60.55095666666667
$ws.Range("H13").Value = 181.65287
$ws.Range("I13").Value = 0.1464834753134679
$ws.Range("J13").Value = 0.1464834753134678
$ws.Range("M13").Value = 1.357726666666667
$ws.Range("N13").Value = 4.073180000000001
$ws.Range("O13").Value = 0.001066367973816652
$ws.Range("P13").Value = 0.001066367973816652
$ws.Range("Q13").Value = 82.21164855851113
$ws.Range("R13").Value = 739.9048370266001
$ws.Range("S13").Value = 0.0001562052867676443
$ws.Range("T13").Value = 0.0001562052867676442

# Row 14
$ws.Range("G14").Value = 52.74960833333333
$ws.Range("H14").Value = 158.248825
$ws.Range("I14").Value = 0.1276106336788006
$ws.Range("J14").Value = 0.1276106336788006
$ws.Range("M14").Value = 1021.934916333333
$ws.Range("N14").Value = 3065.804749
$ws.Range("O14").Value = 0.8026347959846111
$ws.Range("P14").Value = 0.802634795984611
$ws.Range("Q14").Value = 53906.6665787411
$ws.Range("R14").Value = 485159.9992086699
$ws.Range("S14").Value = 0.1024247349282511
$ws.Range("T14").Value = 0.1024247349282511

# Row 15
$ws.Range("G15").Value = 52.74960833333333
$ws.Range("H15").Value = 158.248825
$ws.Range("I15").Value = 0.1276106336788006
$ws.Range("J15").Value = 0.1276106336788006
$ws.Range("O15").Value = 0.04931810976893385
$ws.Range("P15").Value = 0.04931810976893384
$ws.Range("Q15").Value = 3312.309549633137
$ws.Range("R15").Value = 29810.78594669823
$ws.Range("S15").Value = 0.006293515239454296
$ws.Range("T15").Value = 0.006293515239454294

# Row 16
$ws.Range("G16").Value = 52.74960833333333
$ws.Range("H16").Value = 158.248825
$ws.Range("I16").Value = 0.1276106336788006
$ws.Range("J16").Value = 0.1276106336788006
$ws.Range("M16").Value = 187.139577
$ws.Range("N16").Value = 561.418731
$ws.Range("O16").Value = 0.1469807262726385
$ws.Range("P16").Value = 0.1469807262726385
$ws.Range("Q16").Value = 9871.539390415675
$ws.Range("R16").Value = 88843.85451374108
$ws.Range("S16").Value = 0.01875630361822174
$ws.Range("T16").Value = 0.01875630361822174

# Row 17
$ws.Range("G17").Value = 52.74960833333333
$ws.Range("H17").Value = 158.248825
$ws.Range("I17").Value = 0.1276106336788006
$ws.Range("J17").Value = 0.1276106336788006
$ws.Range("M17").Value = 1.357726666666667
$ws.Range("N17").Value = 4.073180000000001
$ws.Range("O17").Value = 0.001066367973816652
$ws.Range("P17").Value = 0.001066367973816652
$ws.Range("Q17").Value = 71.61954989038891
$ws.Range("R17").Value = 644.5759490135001
$ws.Range("S17").Value = 0.0001360798928735216
$ws.Range("T17").Value = 0.0001360798928735216
